# update UG and DG and PPP
#
# 1) Refresh the cached "Date placeholder" field text (used on the slide
#    master + every slide layout) from 7/6/2018 to 4/2/19.
# 2) Nudge a handful of shapes on slide 1 to their new positions.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText($shapes, $newText) {
  for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
      $sh.TextFrame.TextRange.Text = $newText
    }
  }
}

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "4/2/19"

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
  $layout = $master.CustomLayouts.Item($li)
  Set-DatePlaceholderText $layout.Shapes "4/2/19"
}

# --- Reposition shapes on slide 1 ---
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
  $sh = $s.Shapes.Item($i)
  switch ($sh.Name) {
    "TextBox 3" {
      $sh.Left = 11.250039
      $sh.Top  = 336.379181
    }
    "Table 4" {
      $sh.Left = 49.110985
      $sh.Top  = 217.660981
    }
    "Rectangle 6" {
      $sh.Left = 39.267363
      $sh.Top  = 209.47902
    }
    "Straight Arrow Connector 2" {
      $sh.Left = 118.83067
      $sh.Top  = 270.0
    }
  }
}
